$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 ("Task 3.5") was re-run; update Mean/Std/Seed_0..Seed_29 with new values.
$row12Updates = @{
    "B12" = 0.01064461902086103
    "C12" = 0.0005977783134932873
    "D12" = 0.01023751881235292
    "E12" = 0.01069353597888126
    "F12" = 0.01090662505092131
    "G12" = 0.01142638351517632
    "H12" = 0.01042385456887589
    "I12" = 0.01089746446422633
    "J12" = 0.01197185602434062
    "K12" = 0.01040914435577777
    "L12" = 0.01069215975615522
    "M12" = 0.01000714674972099
    "N12" = 0.01015238800555611
    "O12" = 0.01087924733372676
    "P12" = 0.0097993223693
    "Q12" = 0.01007686997159644
    "R12" = 0.01018823998682908
    "S12" = 0.0102420787240653
    "T12" = 0.01042064613234451
    "U12" = 0.01061791013628076
    "V12" = 0.01031125561989872
    "W12" = 0.0101727840990696
    "X12" = 0.01058454807166361
    "Y12" = 0.01065592379424043
    "Z12" = 0.01169337690567628
    "AA12" = 0.01246221382894587
    "AB12" = 0.01024881006312314
    "AC12" = 0.01042735895887869
    "AD12" = 0.01086263077957389
    "AE12" = 0.009964249901595657
    "AF12" = 0.01070290221748107
    "AG12" = 0.01121012444955634
}

foreach ($addr in $row12Updates.Keys) {
    $ws.Range($addr).Value = $row12Updates[$addr]
}

# New row 24 ("Task 4") appended below the existing data.
$row24Updates = @{
    "A24" = "Task 4"
    "B24" = 0.01255425673669903
    "C24" = 0.0006583993878467496
    "D24" = 0.01228864237644526
    "E24" = 0.01239493831392954
    "F24" = 0.01320084748559273
    "G24" = 0.0135264380862166
    "H24" = 0.01178920175774237
    "I24" = 0.01205237801417426
    "J24" = 0.0138306402513558
    "K24" = 0.01230702130655395
    "L24" = 0.0119169713591391
    "M24" = 0.0126006423659272
    "N24" = 0.01220005656477842
    "O24" = 0.0123339049428252
    "P24" = 0.01224519318585674
    "Q24" = 0.01175575738206552
    "R24" = 0.0127559384810645
    "S24" = 0.0132778541071266
    "T24" = 0.01203126105397126
    "U24" = 0.0121738441414118
    "V24" = 0.01186061735295003
    "W24" = 0.0135572255948135
    "X24" = 0.01261200975557703
    "Y24" = 0.01213630682630265
    "Z24" = 0.01201770924054356
    "AA24" = 0.01269238361796149
    "AB24" = 0.0121447554838243
    "AC24" = 0.01258021438784777
    "AD24" = 0.01309439147371852
    "AE24" = 0.01190520251928036
    "AF24" = 0.01449880473366081
    "AG24" = 0.01284654993831405
}

foreach ($addr in $row24Updates.Keys) {
    $ws.Range($addr).Value = $row24Updates[$addr]
}
